# Weekly update: insert two new rows of fresh "Apio" (celery) price data
# at the top of the Vega Central Mapocho de Santiago block (rows 293-294),
# pushing the previously-existing rows 293-303 down to 295-305.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 293 (shifts 293:303 -> 295:305).
$ws.Rows.Item(293).Insert()
$ws.Rows.Item(294).Insert()

# --- New row 293 (Primera) ---
$ws.Cells.Item(293, 1).Value2 = 9
$ws.Cells.Item(293, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(293, 3).Value2 = "Metropolitana"
$ws.Cells.Item(293, 4).Value2 = 44747
$ws.Cells.Item(293, 5).Value2 = 13
$ws.Cells.Item(293, 6).Value2 = 100112017
$ws.Cells.Item(293, 7).Value2 = "Apio"
$ws.Cells.Item(293, 8).Value2 = "Americana (o)"
$ws.Cells.Item(293, 9).Value2 = "Primera"
$ws.Cells.Item(293, 10).Value2 = 70
$ws.Cells.Item(293, 11).Value2 = 8000
$ws.Cells.Item(293, 12).Value2 = 9000
$ws.Cells.Item(293, 13).Value2 = 8500
$ws.Cells.Item(293, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(293, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(293, 16).Value2 = 1417
$ws.Cells.Item(293, 17).Value2 = 6
$ws.Cells.Item(293, 18).Value2 = "Hortaliza"

# --- New row 294 (Segunda) ---
$ws.Cells.Item(294, 1).Value2 = 9
$ws.Cells.Item(294, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(294, 3).Value2 = "Metropolitana"
$ws.Cells.Item(294, 4).Value2 = 44747
$ws.Cells.Item(294, 5).Value2 = 13
$ws.Cells.Item(294, 6).Value2 = 100112017
$ws.Cells.Item(294, 7).Value2 = "Apio"
$ws.Cells.Item(294, 8).Value2 = "Americana (o)"
$ws.Cells.Item(294, 9).Value2 = "Segunda"
$ws.Cells.Item(294, 10).Value2 = 34
$ws.Cells.Item(294, 11).Value2 = 7000
$ws.Cells.Item(294, 12).Value2 = 7000
$ws.Cells.Item(294, 13).Value2 = 7000
$ws.Cells.Item(294, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(294, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(294, 16).Value2 = 1167
$ws.Cells.Item(294, 17).Value2 = 6
$ws.Cells.Item(294, 18).Value2 = "Hortaliza"

Write-Output "done"
